$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 26-35 are a weekly time series for the same market/product.
# Each row's date/volume/price values shift up by one row (row N takes
# the data previously held by row N+1), and a new data point is appended
# at the bottom (row 35).
#
# Columns: D = Fecha, J = Volumen, K = Precio mínimo, L = Precio máximo,
#          M = Precio promedio ponderado, P = Precio $/Kg

$rows = @(
    @{ Row = 26; D = 44260; J = 250; K = 900;  L = 1000; M = 950;  P = 950  }
    @{ Row = 27; D = 44302; J = 200; K = 900;  L = 1000; M = 950;  P = 950  }
    @{ Row = 28; D = 44274; J = 250; K = 1000; L = 1200; M = 1100; P = 1100 }
    @{ Row = 29; D = 44280; J = 250; K = 1400; L = 1500; M = 1450; P = 1450 }
    @{ Row = 30; D = 44432; J = 300; K = 2300; L = 2500; M = 2400; P = 2400 }
    @{ Row = 31; D = 44362; J = 250; K = 2800; L = 3000; M = 2900; P = 2900 }
    @{ Row = 32; D = 44747; J = 250; K = 2000; L = 2500; M = 2250; P = 2250 }
    @{ Row = 33; D = 44603; J = 250; K = 2500; L = 3000; M = 2750; P = 2750 }
    @{ Row = 34; D = 44635; J = 300; K = 1900; L = 2000; M = 1950; P = 1950 }
    @{ Row = 35; D = 44917; J = 400; K = 1500; L = 2000; M = 1750; P = 1750 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("D$n").Value = $r.D
    $ws.Range("J$n").Value = $r.J
    $ws.Range("K$n").Value = $r.K
    $ws.Range("L$n").Value = $r.L
    $ws.Range("M$n").Value = $r.M
    $ws.Range("P$n").Value = $r.P
}
